$d = $word.ActiveDocument

$d.Content.Find.Execute('27÷7=3, 6', $true, $true, $false, $false, $false, $true, 1, $false, '98÷7=14, 0', 2) | Out-Null
$d.Content.Find.Execute('51÷5=10, 1', $true, $true, $false, $false, $false, $true, 1, $false, '73÷2=36, 1', 2) | Out-Null
$d.Content.Find.Execute('39÷4=9, 3', $true, $true, $false, $false, $false, $true, 1, $false, '14÷5=2, 4', 2) | Out-Null
$d.Content.Find.Execute('87÷2=43, 1', $true, $true, $false, $false, $false, $true, 1, $false, '60÷4=15, 0', 2) | Out-Null
$d.Content.Find.Execute('40÷9=4, 4', $true, $true, $false, $false, $false, $true, 1, $false, '80÷2=40, 0', 2) | Out-Null
$d.Content.Find.Execute('73÷5=14, 3', $true, $true, $false, $false, $false, $true, 1, $false, '55÷6=9, 1', 2) | Out-Null
$d.Content.Find.Execute('90÷4=22, 2', $true, $true, $false, $false, $false, $true, 1, $false, '65÷8=8, 1', 2) | Out-Null
$d.Content.Find.Execute('29÷4=7, 1', $true, $true, $false, $false, $false, $true, 1, $false, '45÷3=15, 0', 2) | Out-Null
$d.Content.Find.Execute('51÷3=17, 0', $true, $true, $false, $false, $false, $true, 1, $false, '98÷4=24, 2', 2) | Out-Null
$d.Content.Find.Execute('63÷3=21, 0', $true, $true, $false, $false, $false, $true, 1, $false, '23÷8=2, 7', 2) | Out-Null
$d.Content.Find.Execute('68÷2=34, 0', $true, $true, $false, $false, $false, $true, 1, $false, '10÷4=2, 2', 2) | Out-Null
$d.Content.Find.Execute('88÷2=44, 0', $true, $true, $false, $false, $false, $true, 1, $false, '48÷7=6, 6', 2) | Out-Null
$d.Content.Find.Execute('54÷8=6, 6', $true, $true, $false, $false, $false, $true, 1, $false, '51÷6=8, 3', 2) | Out-Null
$d.Content.Find.Execute('97÷9=10, 7', $true, $true, $false, $false, $false, $true, 1, $false, '81÷9=9, 0', 2) | Out-Null
$d.Content.Find.Execute('34÷7=4, 6', $true, $true, $false, $false, $false, $true, 1, $false, '80÷3=26, 2', 2) | Out-Null
$d.Content.Find.Execute('51÷7=7, 2', $true, $true, $false, $false, $false, $true, 1, $false, '85÷3=28, 1', 2) | Out-Null
$d.Content.Find.Execute('69÷4=17, 1', $true, $true, $false, $false, $false, $true, 1, $false, '18÷5=3, 3', 2) | Out-Null
$d.Content.Find.Execute('50÷8=6, 2', $true, $true, $false, $false, $false, $true, 1, $false, '70÷4=17, 2', 2) | Out-Null
$d.Content.Find.Execute('16÷8=2, 0', $true, $true, $false, $false, $false, $true, 1, $false, '97÷2=48, 1', 2) | Out-Null
$d.Content.Find.Execute('80÷9=8, 8', $true, $true, $false, $false, $false, $true, 1, $false, '55÷2=27, 1', 2) | Out-Null
$d.Content.Find.Execute('70÷9=7, 7', $true, $true, $false, $false, $false, $true, 1, $false, '60÷7=8, 4', 2) | Out-Null
$d.Content.Find.Execute('95÷8=11, 7', $true, $true, $false, $false, $false, $true, 1, $false, '41÷9=4, 5', 2) | Out-Null
$d.Content.Find.Execute('49÷8=6, 1', $true, $true, $false, $false, $false, $true, 1, $false, '84÷5=16, 4', 2) | Out-Null
$d.Content.Find.Execute('75÷7=10, 5', $true, $true, $false, $false, $false, $true, 1, $false, '50÷9=5, 5', 2) | Out-Null
$d.Content.Find.Execute('94÷2=47, 0', $true, $true, $false, $false, $false, $true, 1, $false, '23÷5=4, 3', 2) | Out-Null
